$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9. This shifts the existing milestone rows
# (old row 9 onward) down by one, preserving their content/formatting.
$ws.Rows.Item(9).Insert()

# The row that used to be row 9 is now row 10 and carries the formatting
# we want to reuse for the brand-new row 9 (same column styles). Copy
# that row's formatting (columns A:G only) down onto the freshly
# inserted row 9.
$ws.Range("A10:G10").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new milestone entry describing the Bellman-Ford work.
$ws.Cells.Item(9, 1).Value = $null
$ws.Cells.Item(9, 2).Value = "Not a milestone - Completed implementation of Bellman-Ford algorithm and separated the routing table building part of both the dijkstra function and the Bellman-Ford function to a separate helper function"
$ws.Cells.Item(9, 3).Value2 = 45811
$ws.Cells.Item(9, 4).Value = "find_shortest_paths_bellman_ford`nbuild_routing_table_from_distances"
$ws.Cells.Item(9, 5).Value = "ChatGPT was used to make pseudocode that helped me to learn about the Bellman-Ford algorithm, also helped with basic debugging"
$ws.Cells.Item(9, 6).Value = "test_graph - This function does not fully test the Bellman-Ford function yet, it just does a basic adhoc test to see if the function actually works"
$ws.Cells.Item(9, 7).Value = "The code for both the Dijkstra and Bellman-Ford algorithm could use further refactoring to avoid code duplication. The easiest of which has already been done, which was moving the routing table building to build_routing_table_from_distances"

# Match the author's row height for the new entry.
$ws.Rows.Item(9).RowHeight = 75
